# Generate Report for Handoff
# Update the "Latest Handoff"-related timestamps for the most recently
# handed-off file (27d6d7bd-0eca-4740-a37e-fc21c2e6faf2.md) across the
# Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: column G = "Latest HO Xliff Generate Date"; row 5 is the
# 27d6d7bd-... file ("Ready for handoff").
$overview.Range("G5").Value = "2016-11-03 19:00:59"

# zh-cn sheet: column H = "Latest Handoff Datetime"; row 5 is the
# 27d6d7bd-... file.
$zhcn.Range("H5").Value = "2016-11-03 19:00:46"

# de-de sheet: column H = "Latest Handoff Datetime"; row 5 is the
# 27d6d7bd-... file.
$dede.Range("H5").Value = "2016-11-03 19:00:59"
